$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 54, pushing existing rows 54..136 down to 55..137
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
# Non-numeric/text columns mirror the constant pattern used by every other
# row in this table (Membrillo / Vega Modelo de Temuco / O'Higgins, etc.).
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44645
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100104
$ws.Cells.Item(54, 8).Value = "Frutos de pepita"
$ws.Cells.Item(54, 9).Value = 100104003
$ws.Cells.Item(54, 10).Value = "Membrillo"
$ws.Cells.Item(54, 11).Value = "Champion"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 80
$ws.Cells.Item(54, 14).Value = 14000
$ws.Cells.Item(54, 15).Value = 14000
$ws.Cells.Item(54, 16).Value = 14000
$ws.Cells.Item(54, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(54, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(54, 19).Value = 778
$ws.Cells.Item(54, 20).Value = 18
